$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 40, shifting existing rows 40-45 down to 41-46.
$ws.Rows.Item(40).Insert()

# Populate the newly inserted row 40 with the new weekly record.
$ws.Range("A40").Value2 = 11
$ws.Range("B40").Value2 = "Vega Monumental Concepción"
$ws.Range("C40").Value2 = "Bíobío"
$ws.Range("D40").Value2 = 44504
$ws.Range("D40").NumberFormat = $ws.Range("D41").NumberFormat
$ws.Range("E40").Value2 = 8
$ws.Range("F40").Value2 = 100112001
$ws.Range("G40").Value2 = "Berenjena"
$ws.Range("H40").Value2 = "Sin especificar"
$ws.Range("I40").Value2 = "Primera"
$ws.Range("J40").Value2 = 150
$ws.Range("K40").Value2 = 8500
$ws.Range("L40").Value2 = 9000
$ws.Range("M40").Value2 = 8767
$ws.Range("N40").Value2 = "$/caja 60 unidades"
$ws.Range("O40").Value2 = "Región de Arica y Parinacota"
$ws.Range("P40").Value2 = 146
$ws.Range("Q40").Value2 = 60
$ws.Range("R40").Value2 = "Hortaliza"
